# Add a new "Summary (N = 1197)" results column.
#
# The existing "Summary" column (B) holds the section headers (with an
# empty value) as well as the actual per-row summary statistics. A new
# column C is introduced to hold the N=1197 summary statistics: for
# every row that has an actual (non-blank) value in column B, that
# value is moved into column C and column B is cleared out. Rows that
# are section headers / blank separators (where column B has no real
# value) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 59

# New header for column C on row 1; column B keeps its existing "Summary" header.
$ws.Cells.Item(1, 3).Value2 = "Summary (N = 1197)"

for ($r = 2; $r -le $lastRow; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $bVal = $bCell.Value2

    if ($bVal -ne $null -and $bVal -ne "") {
        $ws.Cells.Item($r, 3).Value2 = $bVal
        $bCell.ClearContents()
    }
}
